# Reran all processing steps after manually fixing wrong recordings date
# information based on app logs. This updates the statistics results
# (normality, equal_var, anova, pairwise_ttests) for the cortisol / weekend
# feature set.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: normality  (Shapiro-Wilk test for each saliva_feature x Weekday/Weekend)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("normality")

$ws.Range("C3").Value = 0.9827
$ws.Range("D3").Value = 0.1106

$ws.Range("C4").Value = 0.9644
$ws.Range("D4").Value = 0.19

$ws.Range("C5").Value = 0.9898
$ws.Range("D5").Value = 0.4837

$ws.Range("C6").Value = 0.975
$ws.Range("D6").Value = 0.448

$ws.Range("C7").Value = 0.9566
$ws.Range("D7").Value = 0.0005

$ws.Range("C8").Value = 0.954
$ws.Range("D8").Value = 0.07729999999999999

$ws.Range("C9").Value = 0.9673
$ws.Range("D9").Value = 0.004

$ws.Range("C10").Value = 0.9651999999999999
$ws.Range("D10").Value = 0.2036

$ws.Range("C11").Value = 0.9778
$ws.Range("D11").Value = 0.037

$ws.Range("C12").Value = 0.9723000000000001
$ws.Range("D12").Value = 0.3637

$ws.Range("C13").Value = 0.9641999999999999
$ws.Range("D13").Value = 0.0022

$ws.Range("C14").Value = 0.9549
$ws.Range("D14").Value = 0.0835

$ws.Range("C15").Value = 0.9529
$ws.Range("D15").Value = 0.0003

$ws.Range("C16").Value = 0.9340000000000001
$ws.Range("D16").Value = 0.0144

# ---------------------------------------------------------------
# Sheet: equal_var  (Levene's test for homoscedasticity)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("equal_var")

$ws.Range("C3").Value = 1.9311
$ws.Range("D3").Value = 0.1665

$ws.Range("C4").Value = 6.2729
$ws.Range("D4").Value = 0.0132

$ws.Range("C5").Value = 4.7792
$ws.Range("D5").Value = 0.0302

$ws.Range("C6").Value = 9.2033
$ws.Range("D6").Value = 0.0028

$ws.Range("C7").Value = 4.0105
$ws.Range("D7").Value = 0.0468
$ws.Range("E7").Value = $false

$ws.Range("C8").Value = 6.7598
$ws.Range("D8").Value = 0.0102

$ws.Range("C9").Value = 6.2213
$ws.Range("D9").Value = 0.0136

# ---------------------------------------------------------------
# Sheet: anova
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("anova")

$ws.Range("E3").Value = 167
$ws.Range("F3").Value = 0.0264
$ws.Range("G3").Value = 0.871

$ws.Range("E4").Value = 167
$ws.Range("F4").Value = 0.0284
$ws.Range("G4").Value = 0.8663999999999999
$ws.Range("H4").Value = 0.0002

$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 0.0043
$ws.Range("G5").Value = 0.9479
$ws.Range("H5").Value = 0

$ws.Range("E6").Value = 167
$ws.Range("F6").Value = 0.6427
$ws.Range("G6").Value = 0.4239
$ws.Range("H6").Value = 0.0038

$ws.Range("E7").Value = 167
$ws.Range("F7").Value = 1.0067
$ws.Range("G7").Value = 0.3172
$ws.Range("H7").Value = 0.006

$ws.Range("E8").Value = 167
$ws.Range("F8").Value = 0.1694
$ws.Range("G8").Value = 0.6812
$ws.Range("H8").Value = 0.001

$ws.Range("E9").Value = 167
$ws.Range("F9").Value = 0.1779
$ws.Range("G9").Value = 0.6737
$ws.Range("H9").Value = 0.0011

# ---------------------------------------------------------------
# Sheet: pairwise_ttests
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("pairwise_ttests")

$ws.Range("H3").Value = 0.1866
$ws.Range("I3").Value = 100.4611
$ws.Range("K3").Value = 0.8524
$ws.Range("L3").Value = "0.19"
$ws.Range("M3").Value = 0.0284
$ws.Range("N3").Value = 0.8524

$ws.Range("H4").Value = 0.201
$ws.Range("I4").Value = 110.3325
$ws.Range("K4").Value = 0.8411
$ws.Range("L4").Value = "0.191"
$ws.Range("M4").Value = 0.0294
$ws.Range("N4").Value = 0.8411

$ws.Range("H5").Value = -0.075
$ws.Range("I5").Value = 100.3238
$ws.Range("K5").Value = 0.9404
$ws.Range("L5").Value = "0.188"
$ws.Range("M5").Value = -0.0114
$ws.Range("N5").Value = 0.9404

$ws.Range("H6").Value = 0.9962
$ws.Range("I6").Value = 122.1799
$ws.Range("K6").Value = 0.3211
$ws.Range("L6").Value = "0.294"
$ws.Range("M6").Value = 0.1399
$ws.Range("N6").Value = 0.3211

$ws.Range("H7").Value = 1.1896
$ws.Range("I7").Value = 108.729
$ws.Range("K7").Value = 0.2368
$ws.Range("L7").Value = "0.356"
$ws.Range("M7").Value = 0.1751
$ws.Range("N7").Value = 0.2368

$ws.Range("H8").Value = 0.502
$ws.Range("I8").Value = 116.6227
$ws.Range("K8").Value = 0.6166
$ws.Range("L8").Value = "0.21"
$ws.Range("M8").Value = 0.0718
$ws.Range("N8").Value = 0.6166

$ws.Range("H9").Value = 0.5034
$ws.Range("I9").Value = 110.5128
$ws.Range("K9").Value = 0.6157
$ws.Range("L9").Value = "0.21"
$ws.Range("M9").Value = 0.0736
$ws.Range("N9").Value = 0.6157
